$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I ("I0") and J ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (border + bold + alignment) from the existing H1
# header cell onto the two new header cells, without touching the values
# we just set (xlPasteFormats = -4122).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Data rows 2-57: new columns I and J ---
$iValues = @(8,5,7,6,4,7,2,6,7,5,8,7,8,5,6,5,9,5,3,8,10,6,6,8,6,9,6,8,9,8,9,3,5,7,8,7,8,9,6,9,7,8,9,8,8,8,9,6,7,6,4,7,4,3,8,4)
$jValues = @(8,6,8,8,4,7,2,7,9,5,8,9,9,6,6,7,9,7,5,8,10,9,7,8,7,9,7,8,9,8,9,5,6,7,8,7,8,9,7,9,7,8,9,8,8,8,9,6,7,6,4,8,4,4,8,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
